$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "27.100.81"; E = "  -1.47%  " },
    @{ Row = 3; D = "1.780.60"; E = "  -2.17%  " },
    @{ Row = 4; D = "1.006"; E = "  +0.16%  " },
    @{ Row = 5; D = "336.34"; E = "  -2.70%  " },
    @{ Row = 6; D = $null; E = "  +0.15%  " },
    @{ Row = 7; D = "0.3812"; E = "  -0.34%  " },
    @{ Row = 8; D = "0.3404"; E = "  -3.42%  " },
    @{ Row = 9; D = "48.12"; E = "  -2.85%  " },
    @{ Row = 10; D = "1.185"; E = "  -4.39%  " },
    @{ Row = 11; D = "0.07422"; E = "  -5.01%  " },
    @{ Row = 12; D = $null; E = "  +0.06%  " },
    @{ Row = 13; D = "21.60"; E = "  -3.06%  " },
    @{ Row = 14; D = "6.429"; E = "  -3.33%  " },
    @{ Row = 15; D = "1.779.43"; E = "  -2.13%  " },
    @{ Row = 16; D = "7.046"; E = "  -3.03%  " },
    @{ Row = 17; D = "0.00001083"; E = "  -4.09%  " },
    @{ Row = 18; D = "0.06631"; E = "  -1.47%  " },
    @{ Row = 19; D = "83.18"; E = "  -3.78%  " },
    @{ Row = 20; D = $null; E = "  +0.15%  " },
    @{ Row = 21; D = "6.532"; E = "  -0.43%  " },
    @{ Row = 22; D = "17.27"; E = "  -2.57%  " },
    @{ Row = 23; D = "27.117.76"; E = "  -1.50%  " },
    @{ Row = 24; D = "12.20"; E = "  -8.51%  " },
    @{ Row = 25; D = "2.377"; E = "  -3.29%  " },
    @{ Row = 26; D = "2.499"; E = "  -7.08%  " },
    @{ Row = 27; D = "1.454"; E = "  -4.03%  " },
    @{ Row = 28; D = "21.02"; E = "  -5.50%  " },
    @{ Row = 29; D = "155.30"; E = "  +1.13%  " },
    @{ Row = 30; D = "1.980.09"; E = "  -2.02%  " },
    @{ Row = 31; D = "133.99"; E = "  -2.20%  " },
    @{ Row = 32; D = "3.988"; E = "  -2.35%  " },
    @{ Row = 33; D = "6.008"; E = "  -6.12%  " },
    @{ Row = 34; D = "0.08662"; E = "  -1.47%  " },
    @{ Row = 35; D = "13.05"; E = "  -7.55%  " },
    @{ Row = 36; D = "1.623"; E = "  -4.53%  " },
    @{ Row = 37; D = "0.6820"; E = "  -3.83%  " },
    @{ Row = 38; D = "5.371"; E = "  -5.44%  " },
    @{ Row = 39; D = "0.06277"; E = "  -4.18%  " },
    @{ Row = 40; D = "0.2170"; E = "  -5.18%  " },
    @{ Row = 41; D = "0.02312"; E = "  -4.69%  " },
    @{ Row = 42; D = "8.535"; E = "  -5.49%  " },
    @{ Row = 43; D = "1.229"; E = "  -5.49%  " },
    @{ Row = 44; D = "14.16"; E = "  -4.51%  " },
    @{ Row = 45; D = $null; E = "  +0.13%  " },
    @{ Row = 46; D = "0.6410"; E = "  -3.36%  " },
    @{ Row = 47; D = "3.856"; E = "  -3.08%  " },
    @{ Row = 48; D = "2.116"; E = "  -3.63%  " },
    @{ Row = 49; D = "131.03"; E = "  -1.50%  " },
    @{ Row = 50; D = "0.07086"; E = "  -3.95%  " },
    @{ Row = 51; D = "78.46"; E = "  -3.10%  " }

)

foreach ($u in $updates) {
    $row = $u.Row
    if ($null -ne $u.D) {
        $ws.Cells.Item($row, 4).Value = $u.D
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($row, 5).Value = $u.E
    }
}
